{"js": "// Replace \"Ich brauche ein Dokument ich ver\u00e4ndere was\" with\n// \"Ich will hier doch noch was anderes haben\", keeping the _GoBack\n// bookmark in place (it originally sat between \" ver\u00e4ndere\" and \" was\").\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) Remove the trailing \" was\" run that sits AFTER the bookmark, so the\n//    bookmark remains the last thing in the paragraph (matches the target).\nconst tail = body.search(\" was\", { matchCase: true, matchWholeWord: false });\ntail.load(\"items\");\nawait context.sync();\n\nfor (const r of tail.items) {\n  r.delete();\n}\nawait context.sync();\n\n// 2) Replace everything from the start of the paragraph up to (but not\n//    including) the bookmark with the new sentence. Using expandTo keeps the\n//    bookmark itself untouched/unmoved in the run structure.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nconst paragraphStart = firstParagraph.getRange(\"Start\");\n\nconst bookmarkRange = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\nbookmarkRange.load(\"isNullObject\");\nawait context.sync();\n\nconst replaceRange = bookmarkRange.isNullObject\n  ? firstParagraph.getRange(\"Whole\")\n  : paragraphStart.expandTo(bookmarkRange);\n\nreplaceRange.insertText(\"Ich will hier doch noch was anderes haben\", \"Replace\");\nawait context.sync();\n", "ps1": "# Replace \"Ich brauche ein Dokument ich ver\u00e4ndere was\" with\n# \"Ich will hier doch noch was anderes haben\", while leaving the\n# _GoBack bookmark sitting exactly where it was (right after \"ver\u00e4ndere\",\n# i.e. now right after the new sentence, with nothing following it).\n\n$d = $word.ActiveDocument\n\n# 1) First remove the trailing \" was\" run that lives AFTER the bookmark, so\n#    it doesn't also get swept up by (or collide with) the main replace\n#    below, and so the bookmark ends up as the last thing in the paragraph.\n$tailRange = $d.Content\n$tailRange.Find.Execute(\" was\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# 2) Replace the remaining original text (everything before the bookmark)\n#    with the new sentence. This leaves the bookmark itself untouched.\n$mainRange = $d.Content\n$mainRange.Find.Execute(\"Ich brauche ein Dokument ich ver\u00e4ndere\", $false, $false, $false, $false, $false, $true, 1, $false, \"Ich will hier doch noch was anderes haben\", 2)\n"}
